$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17:A21").Value = "Resolving-Mac"

$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.07057933333333334
$ws.Cells.Item(2, 8).Value = 0.211738
$ws.Cells.Item(2, 9).Value = 0.1281663575568867
$ws.Cells.Item(2, 10).Value = 0.1281663575568867
$ws.Cells.Item(2, 13).Value = 3.770298333333333
$ws.Cells.Item(2, 14).Value = 11.310895
$ws.Cells.Item(2, 15).Value = 0.06100259562224731
$ws.Cells.Item(2, 16).Value = 0.06125631726190612
$ws.Cells.Item(2, 17).Value = 0.2661051428344445
$ws.Cells.Item(2, 18).Value = 2.39494628551
$ws.Cells.Item(2, 19).Value = 0.007818480482419121
$ws.Cells.Item(2, 20).Value = 0.007850999060807551
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.07057933333333334
$ws.Cells.Item(3, 8).Value = 0.211738
$ws.Cells.Item(3, 9).Value = 0.1281663575568867
$ws.Cells.Item(3, 10).Value = 0.1281663575568867
$ws.Cells.Item(3, 15).Value = 0.06469423882843597
$ws.Cells.Item(3, 16).Value = 0.06496331472897099
$ws.Cells.Item(3, 17).Value = 0.2822088058451112
$ws.Cells.Item(3, 18).Value = 2.539879252606
$ws.Cells.Item(3, 19).Value = 0.008291624945555949
$ws.Cells.Item(3, 20).Value = 0.008326111423633861
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.07057933333333334
$ws.Cells.Item(4, 8).Value = 0.211738
$ws.Cells.Item(4, 9).Value = 0.1281663575568867
$ws.Cells.Item(4, 10).Value = 0.1281663575568867
$ws.Cells.Item(4, 13).Value = 31.40746233333333
$ws.Cells.Item(4, 14).Value = 94.222387
$ws.Cells.Item(4, 15).Value = 0.5081658147055464
$ws.Cells.Item(4, 16).Value = 0.5102793749960634
$ws.Cells.Item(4, 17).Value = 2.216717753178445
$ws.Cells.Item(4, 18).Value = 19.950459778606
$ws.Cells.Item(4, 19).Value = 0.06512976150573771
$ws.Cells.Item(4, 20).Value = 0.06540064882965013
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.07057933333333334
$ws.Cells.Item(5, 8).Value = 0.211738
$ws.Cells.Item(5, 9).Value = 0.1281663575568867
$ws.Cells.Item(5, 10).Value = 0.1281663575568867
$ws.Cells.Item(5, 13).Value = 0.7679895
$ws.Cells.Item(5, 14).Value = 1.535979
$ws.Cells.Item(5, 15).Value = 0.01242590075603175
$ws.Cells.Item(5, 16).Value = 0.008318388326620067
$ws.Cells.Item(5, 17).Value = 0.05420418691700001
$ws.Cells.Item(5, 18).Value = 0.325225121502
$ws.Cells.Item(5, 19).Value = 0.001592582439263954
$ws.Cells.Item(5, 20).Value = 0.00106613753256662
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.07057933333333334
$ws.Cells.Item(6, 8).Value = 0.211738
$ws.Cells.Item(6, 9).Value = 0.1281663575568867
$ws.Cells.Item(6, 10).Value = 0.1281663575568867
$ws.Cells.Item(6, 13).Value = 21.861327
$ws.Cells.Item(6, 14).Value = 65.58398100000001
$ws.Cells.Item(6, 15).Value = 0.3537114500877385
$ws.Cells.Item(6, 16).Value = 0.3551826046864394
$ws.Cells.Item(6, 17).Value = 1.542957885442
$ws.Cells.Item(6, 18).Value = 13.886620968978
$ws.Cells.Item(6, 19).Value = 0.04533390818390998
$ws.Cells.Item(6, 20).Value = 0.04552246071022854
$ws.Cells.Item(7, 9).Value = 0.05066414213561767
$ws.Cells.Item(7, 10).Value = 0.05066414213561768
$ws.Cells.Item(7, 13).Value = 3.770298333333333
$ws.Cells.Item(7, 14).Value = 11.310895
$ws.Cells.Item(7, 15).Value = 0.06100259562224731
$ws.Cells.Item(7, 16).Value = 0.06125631726190612
$ws.Cells.Item(7, 17).Value = 0.1051913235
$ws.Cells.Item(7, 18).Value = 0.9467219115
$ws.Cells.Item(7, 19).Value = 0.003090644175247146
$ws.Cells.Item(7, 20).Value = 0.003103498764461702
$ws.Cells.Item(8, 9).Value = 0.05066414213561767
$ws.Cells.Item(8, 10).Value = 0.05066414213561768
$ws.Cells.Item(8, 15).Value = 0.06469423882843597
$ws.Cells.Item(8, 16).Value = 0.06496331472897099
$ws.Cells.Item(8, 19).Value = 0.003277678111359476
$ws.Cells.Item(8, 20).Value = 0.003291310611029452
$ws.Cells.Item(9, 9).Value = 0.05066414213561767
$ws.Cells.Item(9, 10).Value = 0.05066414213561768
$ws.Cells.Item(9, 13).Value = 31.40746233333333
$ws.Cells.Item(9, 14).Value = 94.222387
$ws.Cells.Item(9, 15).Value = 0.5081658147055464
$ws.Cells.Item(9, 16).Value = 0.5102793749960634
$ws.Cells.Item(9, 17).Value = 0.8762681990999999
$ws.Cells.Item(9, 18).Value = 7.8864137919
$ws.Cells.Item(9, 19).Value = 0.02574578506470376
$ws.Cells.Item(9, 20).Value = 0.02585286678367471
$ws.Cells.Item(10, 9).Value = 0.05066414213561767
$ws.Cells.Item(10, 10).Value = 0.05066414213561768
$ws.Cells.Item(10, 13).Value = 0.7679895
$ws.Cells.Item(10, 14).Value = 1.535979
$ws.Cells.Item(10, 15).Value = 0.01242590075603175
$ws.Cells.Item(10, 16).Value = 0.008318388326620067
$ws.Cells.Item(10, 17).Value = 0.02142690705
$ws.Cells.Item(10, 18).Value = 0.1285614423
$ws.Cells.Item(10, 19).Value = 0.0006295476020666715
$ws.Cells.Item(10, 20).Value = 0.000421444008519142
$ws.Cells.Item(11, 9).Value = 0.05066414213561767
$ws.Cells.Item(11, 10).Value = 0.05066414213561768
$ws.Cells.Item(11, 13).Value = 21.861327
$ws.Cells.Item(11, 14).Value = 65.58398100000001
$ws.Cells.Item(11, 15).Value = 0.3537114500877385
$ws.Cells.Item(11, 16).Value = 0.3551826046864394
$ws.Cells.Item(11, 17).Value = 0.6099310233
$ws.Cells.Item(11, 18).Value = 5.489379209700001
$ws.Cells.Item(11, 19).Value = 0.01792048718224062
$ws.Cells.Item(11, 20).Value = 0.01799502196793267
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.128589
$ws.Cells.Item(12, 8).Value = 0.385767
$ws.Cells.Item(12, 9).Value = 0.2335072176730087
$ws.Cells.Item(12, 10).Value = 0.2335072176730087
$ws.Cells.Item(12, 13).Value = 3.770298333333333
$ws.Cells.Item(12, 14).Value = 11.310895
$ws.Cells.Item(12, 15).Value = 0.06100259562224731
$ws.Cells.Item(12, 16).Value = 0.06125631726190612
$ws.Cells.Item(12, 17).Value = 0.484818892385
$ws.Cells.Item(12, 18).Value = 4.363370031465
$ws.Cells.Item(12, 19).Value = 0.01424454637458263
$ws.Cells.Item(12, 20).Value = 0.01430379220872279
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.128589
$ws.Cells.Item(13, 8).Value = 0.385767
$ws.Cells.Item(13, 9).Value = 0.2335072176730087
$ws.Cells.Item(13, 10).Value = 0.2335072176730087
$ws.Cells.Item(13, 15).Value = 0.06469423882843597
$ws.Cells.Item(13, 16).Value = 0.06496331472897099
$ws.Cells.Item(13, 17).Value = 0.5141582729810001
$ws.Cells.Item(13, 18).Value = 4.627424456829001
$ws.Cells.Item(13, 19).Value = 0.01510657170830121
$ws.Cells.Item(13, 20).Value = 0.015169402873178
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.128589
$ws.Cells.Item(14, 8).Value = 0.385767
$ws.Cells.Item(14, 9).Value = 0.2335072176730087
$ws.Cells.Item(14, 10).Value = 0.2335072176730087
$ws.Cells.Item(14, 13).Value = 31.40746233333333
$ws.Cells.Item(14, 14).Value = 94.222387
$ws.Cells.Item(14, 15).Value = 0.5081658147055464
$ws.Cells.Item(14, 16).Value = 0.5102793749960634
$ws.Cells.Item(14, 17).Value = 4.038654173981
$ws.Cells.Item(14, 18).Value = 36.34788756582901
$ws.Cells.Item(14, 19).Value = 0.1186603855084298
$ws.Cells.Item(14, 20).Value = 0.1191539170912526
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.128589
$ws.Cells.Item(15, 8).Value = 0.385767
$ws.Cells.Item(15, 9).Value = 0.2335072176730087
$ws.Cells.Item(15, 10).Value = 0.2335072176730087
$ws.Cells.Item(15, 13).Value = 0.7679895
$ws.Cells.Item(15, 14).Value = 1.535979
$ws.Cells.Item(15, 15).Value = 0.01242590075603175
$ws.Cells.Item(15, 16).Value = 0.008318388326620067
$ws.Cells.Item(15, 17).Value = 0.09875500181550001
$ws.Cells.Item(15, 18).Value = 0.5925300108930001
$ws.Cells.Item(15, 19).Value = 0.002901537512621908
$ws.Cells.Item(15, 20).Value = 0.001942403713672687
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.128589
$ws.Cells.Item(16, 8).Value = 0.385767
$ws.Cells.Item(16, 9).Value = 0.2335072176730087
$ws.Cells.Item(16, 10).Value = 0.2335072176730087
$ws.Cells.Item(16, 13).Value = 21.861327
$ws.Cells.Item(16, 14).Value = 65.58398100000001
$ws.Cells.Item(16, 15).Value = 0.3537114500877385
$ws.Cells.Item(16, 16).Value = 0.3551826046864394
$ws.Cells.Item(16, 17).Value = 2.811126177603001
$ws.Cells.Item(16, 18).Value = 25.300135598427
$ws.Cells.Item(16, 19).Value = 0.0825941765690731
$ws.Cells.Item(16, 20).Value = 0.08293770178618259
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.323617
$ws.Cells.Item(17, 8).Value = 0.970851
$ws.Cells.Item(17, 9).Value = 0.5876622826344869
$ws.Cells.Item(17, 10).Value = 0.5876622826344869
$ws.Cells.Item(17, 13).Value = 3.770298333333333
$ws.Cells.Item(17, 14).Value = 11.310895
$ws.Cells.Item(17, 15).Value = 0.06100259562224731
$ws.Cells.Item(17, 16).Value = 0.06125631726190612
$ws.Cells.Item(17, 17).Value = 1.220132635738333
$ws.Cells.Item(17, 18).Value = 10.981193721645
$ws.Cells.Item(17, 19).Value = 0.03584892458999842
$ws.Cells.Item(17, 20).Value = 0.03599802722791408
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0.323617
$ws.Cells.Item(18, 8).Value = 0.970851
$ws.Cells.Item(18, 9).Value = 0.5876622826344869
$ws.Cells.Item(18, 10).Value = 0.5876622826344869
$ws.Cells.Item(18, 15).Value = 0.06469423882843597
$ws.Cells.Item(18, 16).Value = 0.06496331472897099
$ws.Cells.Item(18, 17).Value = 1.293970384926334
$ws.Cells.Item(18, 18).Value = 11.645733464337
$ws.Cells.Item(18, 19).Value = 0.03801836406321934
$ws.Cells.Item(18, 20).Value = 0.03817648982112968
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 0.323617
$ws.Cells.Item(19, 8).Value = 0.970851
$ws.Cells.Item(19, 9).Value = 0.5876622826344869
$ws.Cells.Item(19, 10).Value = 0.5876622826344869
$ws.Cells.Item(19, 13).Value = 31.40746233333333
$ws.Cells.Item(19, 14).Value = 94.222387
$ws.Cells.Item(19, 15).Value = 0.5081658147055464
$ws.Cells.Item(19, 16).Value = 0.5102793749960634
$ws.Cells.Item(19, 17).Value = 10.16398873792633
$ws.Cells.Item(19, 18).Value = 91.475898641337
$ws.Cells.Item(19, 19).Value = 0.2986298826266752
$ws.Cells.Item(19, 20).Value = 0.2998719422914859
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 0.323617
$ws.Cells.Item(20, 8).Value = 0.970851
$ws.Cells.Item(20, 9).Value = 0.5876622826344869
$ws.Cells.Item(20, 10).Value = 0.5876622826344869
$ws.Cells.Item(20, 13).Value = 0.7679895
$ws.Cells.Item(20, 14).Value = 1.535979
$ws.Cells.Item(20, 15).Value = 0.01242590075603175
$ws.Cells.Item(20, 16).Value = 0.008318388326620067
$ws.Cells.Item(20, 17).Value = 0.2485344580215
$ws.Cells.Item(20, 18).Value = 1.491206748129
$ws.Cells.Item(20, 19).Value = 0.007302233202079213
$ws.Cells.Item(20, 20).Value = 0.004888403071861619
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 0.323617
$ws.Cells.Item(21, 8).Value = 0.970851
$ws.Cells.Item(21, 9).Value = 0.5876622826344869
$ws.Cells.Item(21, 10).Value = 0.5876622826344869
$ws.Cells.Item(21, 13).Value = 21.861327
$ws.Cells.Item(21, 14).Value = 65.58398100000001
$ws.Cells.Item(21, 15).Value = 0.3537114500877385
$ws.Cells.Item(21, 16).Value = 0.3551826046864394
$ws.Cells.Item(21, 17).Value = 7.074697059759001
$ws.Cells.Item(21, 18).Value = 63.67227353783101
$ws.Cells.Item(21, 19).Value = 0.2078628781525148
$ws.Cells.Item(21, 20).Value = 0.2087274202220956
